$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A5").Value = "Manitoba"
$ws.Range("A6").Value = "Ikeja"
$ws.Range("A7").Value = "Istanbul"
$ws.Range("A8").Value = "Abuja"

$ws.Range("A8").Select()
